# Auto-generated edit script applying numeric corrections from the
# Zodiark_Profits scheduled-runner update across the ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR sheets (profit recalculations for specific leve rows).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3917.2222
$ws.Range("I62").Value = 3700.7334
$ws.Range("K62").Value = 3700.7334
$ws.Range("M62").Value = -3076.7334
$ws.Range("H65").Value = 3917.2222
$ws.Range("I65").Value = 3700.7334
$ws.Range("K65").Value = 18503.667
$ws.Range("M65").Value = -15383.667
$ws.Range("H86").Value = 34003
$ws.Range("I86").Value = 1004.5
$ws.Range("K86").Value = 1004.5
$ws.Range("M86").Value = 118.5
$ws.Range("H89").Value = 34003
$ws.Range("I89").Value = 1004.5
$ws.Range("K89").Value = 5022.5
$ws.Range("M89").Value = 593.5
$ws.Range("H106").Value = 12276.914
$ws.Range("I106").Value = 13593.742
$ws.Range("K106").Value = 13593.742
$ws.Range("M106").Value = -12962.742
$ws.Range("H137").Value = 950.1667
$ws.Range("I137").Value = 893.1
$ws.Range("J137").Value = 1021.5
$ws.Range("K137").Value = 2679.3
$ws.Range("L137").Value = 3064.5
$ws.Range("M137").Value = -129.3000000000002
$ws.Range("N137").Value = -8164.5
$ws.Range("H138").Value = 5444.3125
$ws.Range("J138").Value = 7138.4194
$ws.Range("L138").Value = 21415.2582
$ws.Range("N138").Value = -31695.2582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3710.375
$ws.Range("I32").Value = 2252.8
$ws.Range("J32").Value = 10998.25
$ws.Range("K32").Value = 2252.8
$ws.Range("L32").Value = 10998.25
$ws.Range("M32").Value = -1965.8
$ws.Range("N32").Value = -11572.25
$ws.Range("H61").Value = 3578.8845
$ws.Range("I61").Value = 3262.6924
$ws.Range("K61").Value = 3262.6924
$ws.Range("M61").Value = -3050.6924
$ws.Range("H88").Value = 1608.25
$ws.Range("I88").Value = 1167.8334
$ws.Range("J88").Value = 1872.5
$ws.Range("K88").Value = 1167.8334
$ws.Range("L88").Value = 1872.5
$ws.Range("M88").Value = -761.8334
$ws.Range("N88").Value = -2684.5
$ws.Range("H91").Value = 1608.25
$ws.Range("I91").Value = 1167.8334
$ws.Range("J91").Value = 1872.5
$ws.Range("K91").Value = 1167.8334
$ws.Range("L91").Value = 1872.5
$ws.Range("M91").Value = 236.1666
$ws.Range("N91").Value = -4680.5
$ws.Range("H97").Value = 1037.909
$ws.Range("I97").Value = 774.4286
$ws.Range("K97").Value = 774.4286
$ws.Range("M97").Value = -278.4286
$ws.Range("H122").Value = 6436.4653
$ws.Range("I122").Value = 6377.875
$ws.Range("K122").Value = 19133.625
$ws.Range("M122").Value = -16683.625
$ws.Range("H132").Value = 7440.5605
$ws.Range("I132").Value = 7711.9644
$ws.Range("K132").Value = 23135.8932
$ws.Range("M132").Value = -20605.8932
$ws.Range("H136").Value = 3578.8845
$ws.Range("I136").Value = 3262.6924
$ws.Range("K136").Value = 9788.0772
$ws.Range("M136").Value = -7238.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2117.84
$ws.Range("I134").Value = 2130.95
$ws.Range("K134").Value = 6392.849999999999
$ws.Range("M134").Value = -3857.849999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 886.04877
$ws.Range("I107").Value = 403.23077
$ws.Range("K107").Value = 403.23077
$ws.Range("M107").Value = 1516.76923
$ws.Range("H122").Value = 928.55
$ws.Range("I122").Value = 881.6667
$ws.Range("K122").Value = 2645.0001
$ws.Range("M122").Value = -195.0001000000002
$ws.Range("H132").Value = 1989.0889
$ws.Range("I132").Value = 1772.8918
$ws.Range("K132").Value = 5318.6754
$ws.Range("M132").Value = -2788.6754
$ws.Range("H134").Value = 1965.675
$ws.Range("I134").Value = 1741.4688
$ws.Range("K134").Value = 5224.4064
$ws.Range("M134").Value = -2689.4064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 376.14285
$ws.Range("I92").Value = 237
$ws.Range("K92").Value = 711
$ws.Range("M92").Value = 537
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11975.143
$ws.Range("I132").Value = 12766.5
$ws.Range("K132").Value = 38299.5
$ws.Range("M132").Value = -35769.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1082.8
$ws.Range("I46").Value = 993
$ws.Range("J46").Value = 1105.25
$ws.Range("K46").Value = 993
$ws.Range("L46").Value = 1105.25
$ws.Range("M46").Value = -805
$ws.Range("N46").Value = -1481.25
$ws.Range("H136").Value = 4333.1665
$ws.Range("I136").Value = 2997.25
$ws.Range("J136").Value = 7005
$ws.Range("K136").Value = 8991.75
$ws.Range("L136").Value = 21015
$ws.Range("M136").Value = -6441.75
$ws.Range("N136").Value = -26115

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3724.2
$ws.Range("I81").Value = 2256.5
$ws.Range("J81").Value = 9595
$ws.Range("K81").Value = 4513
$ws.Range("L81").Value = 19190
$ws.Range("M81").Value = -3452
$ws.Range("N81").Value = -21312
$ws.Range("H84").Value = 3724.2
$ws.Range("I84").Value = 2256.5
$ws.Range("J84").Value = 9595
$ws.Range("K84").Value = 22565
$ws.Range("L84").Value = 95950
$ws.Range("M84").Value = -17261
$ws.Range("N84").Value = -106558
$ws.Range("H96").Value = 3089741.5
$ws.Range("I96").Value = 5293418.5
$ws.Range("K96").Value = 5293418.5
$ws.Range("M96").Value = -5292045.5
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H113").Value = 1494.95
$ws.Range("I113").Value = 1505.7646
$ws.Range("K113").Value = 4517.293799999999
$ws.Range("M113").Value = -2347.293799999999
$ws.Range("H122").Value = 3582.0833
$ws.Range("I122").Value = 3473.6453
$ws.Range("K122").Value = 10420.9359
$ws.Range("M122").Value = -7970.9359
$ws.Range("H132").Value = 2883.4666
$ws.Range("I132").Value = 1975.2
$ws.Range("K132").Value = 5925.6
$ws.Range("M132").Value = -3395.6
$ws.Range("H136").Value = 7314.5713
$ws.Range("I136").Value = 3534.5
$ws.Range("K136").Value = 10603.5
$ws.Range("M136").Value = -8053.5
$ws.Range("H140").Value = 100000.125
$ws.Range("J140").Value = 100000.125
$ws.Range("L140").Value = 100000.125
$ws.Range("N140").Value = -110360.125
$ws.Range("H141").Value = 250000
$ws.Range("J141").Value = 250000
$ws.Range("L141").Value = 250000
$ws.Range("N141").Value = -260360
